# feat: add 2022-Q4 data
#
# Workbook currently has two sheets: "总计" (totals) and "2022-Q1" (fund
# holdings detail for 2022-Q1). We are adding a new quarter (2022-Q4):
#   1. The "总计" sheet gets a new top data row for 2022-Q4, and the old
#      2022-Q1 row is pushed down to row 3 (unchanged values).
#   2. The existing "2022-Q1" detail sheet is cloned (so the old Q1 detail
#      data is preserved verbatim in its own "2022-Q1" tab after the move),
#      then the original tab is renamed to "2022-Q4" and its content is
#      replaced with the new quarter's fund holdings.

$wb = $excel.ActiveWorkbook

$totals = $wb.Worksheets.Item(1)
$detail = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert the new 2022-Q4 row, push old 2022-Q1 row down.
# ---------------------------------------------------------------------

# Copy the current row 2 (2022-Q1, count=2, value=0.1) down to row 3 first,
# preserving its formatting (style index carried on column A).
$totals.Range("A3").Value2 = 1
$totals.Range("B3").Value2 = "2022-Q1"
$totals.Range("C3").Value2 = 2
$totals.Range("D3").Value2 = 0.1

$totals.Range("A2").Copy()
$totals.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Now overwrite row 2 with the new 2022-Q4 figures (index/count unchanged).
$totals.Range("B2").Value2 = "2022-Q4"
$totals.Range("D2").Value2 = 0.16

# ---------------------------------------------------------------------
# 2) Duplicate the existing "2022-Q1" detail sheet so its data survives
#    under its own tab once the original tab becomes "2022-Q4".
# ---------------------------------------------------------------------

$detail.Copy([System.Reflection.Missing]::Value, $detail)
$oldQ1 = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 3) Rename the original detail sheet to "2022-Q4" and replace its data
#    with the new quarter's fund holdings.
# ---------------------------------------------------------------------

$detail.Name = "2022-Q4"
$oldQ1.Name = "2022-Q1"

# Clear out the old 2022-Q1 rows (2 and 3) on the now-"2022-Q4" sheet; the
# header row (row 1) and its formatting stay as-is.
$detail.Range("A2:H3").ClearContents()

$detail.Range("B1").Value2 = "基金代码"
$detail.Range("C1").Value2 = "基金名称"
$detail.Range("D1").Value2 = "基金规模"
$detail.Range("E1").Value2 = "股票总仓位"
$detail.Range("F1").Value2 = "仓位占比"
$detail.Range("G1").Value2 = "持有市值(亿元)"
$detail.Range("H1").Value2 = "仓位排名"

$detail.Range("A2").Value2 = 0
$detail.Range("B2").Value2 = "015697"
$detail.Range("C2").Value2 = "华夏磐润两年定开混合A"
$detail.Range("D2").Value2 = "2.68"
$detail.Range("E2").Value2 = "86.76"
$detail.Range("F2").Value2 = "4.37"
$detail.Range("G2").Value2 = "0.1171"
$detail.Range("H2").Value2 = 5

$detail.Range("A3").Value2 = 1
$detail.Range("B3").Value2 = "015698"
$detail.Range("C3").Value2 = "华夏磐润两年定开混合C"
$detail.Range("D3").Value2 = "0.99"
$detail.Range("E3").Value2 = "86.76"
$detail.Range("F3").Value2 = "4.37"
$detail.Range("G3").Value2 = "0.0433"
$detail.Range("H3").Value2 = 5

# Match formatting used for this sheet's header/index cells in the source
# workbook (bold + bordered + centered "header" style, shared with the
# "总计" sheet rather than the style the old detail sheet used).
$totals.Range("B1:D1").Copy()
$detail.Range("B1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totals.Range("B1").Copy()
$detail.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totals.Range("A2").Copy()
$detail.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Page margins on the new "2022-Q4" sheet match the "总计" sheet's margins
# (0.75in/0.75in/1in/1in/0.5in/0.5in) rather than the detail sheet's
# original Excel-default margins.
$detail.PageSetup.LeftMargin = 54
$detail.PageSetup.RightMargin = 54
$detail.PageSetup.TopMargin = 72
$detail.PageSetup.BottomMargin = 72
$detail.PageSetup.HeaderMargin = 36
$detail.PageSetup.FooterMargin = 36

$totals.Select()
